$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header cells for the new columns, matching the style of the existing
# header row (row 1, style index 1 - bold / centered / bordered).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) onto the
# new header cells so they look consistent with the rest of row 1.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the team record (Wins/Losses/Ties) for every data row, 2 through 57.
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 30).Value = 71  # AD
    $ws.Cells.Item($r, 31).Value = 91  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
